# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet, and the
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" entries on
# the per-language sheets for the dc6c0c59-... file row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the dc6c0c59-... file; column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G3").Value = "2016-08-28 20:47:50"

# zh-cn sheet: row 3 is the dc6c0c59-... file
# Column H = "Correspond Handoff Datetime"
$wsZhCn.Range("H3").Value = "2016-08-28 20:47:45"
# Column K = "Correspond Handback DateTime"
$wsZhCn.Range("K3").Value = "2016-08-28 20:48:06"

# de-de sheet: row 3 is the dc6c0c59-... file
# Column H = "Correspond Handoff Datetime"
$wsDeDe.Range("H3").Value = "2016-08-28 20:47:50"
# Column K = "Correspond Handback DateTime"
$wsDeDe.Range("K3").Value = "2016-08-28 20:48:13"
